# The document has a single table. Several rows contain a "Đáp ứng"
# cell, but the one to clear is the vertically-merged ("restart") cell
# in the row for "Phân hệ Chức năng chung của hệ thống" — table row 26,
# column 4 (1-based Word indices; 0-based row 25 / cell 3).
#
# The target change removes the whole <w:r> run (and its <w:t>) from
# that cell's paragraph, leaving an empty paragraph behind — not just
# an empty <w:t/>. So we delete the run's character range explicitly
# (excluding the trailing paragraph mark) rather than assigning "" to
# Range.Text, which would leave a stray empty run.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$cell = $t.Cell(26, 4)

$cellRange = $cell.Range
# Exclude the final paragraph-mark / end-of-cell character so only the
# run's text is targeted, then delete it — this removes the run
# entirely rather than leaving an empty <w:r><w:t/></w:r>.
$textRange = $d.Range($cellRange.Start, $cellRange.End - 1)

if ($textRange.Text -eq "Đáp ứng") {
    $textRange.Delete()
}
